$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.4953129999999999
$ws.Range("H2").Value = 1.485939
$ws.Range("I2").Value = 0.2134181161789063
$ws.Range("J2").Value = 0.2134181161789063
$ws.Range("M2").Value = 22.91402233333334
$ws.Range("N2").Value = 68.74206700000001
$ws.Range("O2").Value = 0.191813973987922
$ws.Range("P2").Value = 0.191813973987922
$ws.Range("Q2").Value = 11.34961314399033
$ws.Range("R2").Value = 102.146518295913
$ws.Range("S2").Value = 0.04093657698529204
$ws.Range("T2").Value = 0.04093657698529205
$ws.Range("G3").Value = 0.4953129999999999
$ws.Range("H3").Value = 1.485939
$ws.Range("I3").Value = 0.2134181161789063
$ws.Range("J3").Value = 0.2134181161789063
$ws.Range("O3").Value = 0.07776078244711707
$ws.Range("P3").Value = 0.07776078244711705
$ws.Range("Q3").Value = 4.601097512344666
$ws.Range("R3").Value = 41.409877611102
$ws.Range("S3").Value = 0.01659555970246148
$ws.Range("T3").Value = 0.01659555970246148
$ws.Range("G4").Value = 0.4953129999999999
$ws.Range("H4").Value = 1.485939
$ws.Range("I4").Value = 0.2134181161789063
$ws.Range("J4").Value = 0.2134181161789063
$ws.Range("M4").Value = 6.020714333333333
$ws.Range("N4").Value = 18.062143
$ws.Range("O4").Value = 0.05039958178109668
$ws.Range("P4").Value = 0.05039958178109668
$ws.Range("Q4").Value = 2.982138078586333
$ws.Range("R4").Value = 26.839242707277
$ws.Range("S4").Value = 0.01075618379992638
$ws.Range("T4").Value = 0.01075618379992638
$ws.Range("G5").Value = 0.4953129999999999
$ws.Range("H5").Value = 1.485939
$ws.Range("I5").Value = 0.2134181161789063
$ws.Range("J5").Value = 0.2134181161789063
$ws.Range("M5").Value = 81.23559966666666
$ws.Range("N5").Value = 243.706799
$ws.Range("O5").Value = 0.6800256617838641
$ws.Range("P5").Value = 0.6800256617838643
$ws.Range("Q5").Value = 40.23704857769566
$ws.Range("R5").Value = 362.133437199261
$ws.Range("S5").Value = 0.1451297956912263
$ws.Range("T5").Value = 0.1451297956912264
$ws.Range("I6").Value = 0.3878984662564351
$ws.Range("J6").Value = 0.3878984662564351
$ws.Range("M6").Value = 22.91402233333334
$ws.Range("N6").Value = 68.74206700000001
$ws.Range("O6").Value = 0.191813973987922
$ws.Range("P6").Value = 0.191813973987922
$ws.Range("Q6").Value = 20.62850900373967
$ws.Range("R6").Value = 185.656581033657
$ws.Range("S6").Value = 0.07440434631646668
$ws.Range("T6").Value = 0.0744043463164667
$ws.Range("I7").Value = 0.3878984662564351
$ws.Range("J7").Value = 0.3878984662564351
$ws.Range("O7").Value = 0.07776078244711707
$ws.Range("P7").Value = 0.07776078244711705
$ws.Range("S7").Value = 0.03016328824613703
$ws.Range("T7").Value = 0.03016328824613703
$ws.Range("I8").Value = 0.3878984662564351
$ws.Range("J8").Value = 0.3878984662564351
$ws.Range("M8").Value = 6.020714333333333
$ws.Range("N8").Value = 18.062143
$ws.Range("O8").Value = 0.05039958178109668
$ws.Range("P8").Value = 0.05039958178109668
$ws.Range("Q8").Value = 5.420190223583666
$ws.Range("R8").Value = 48.781712012253
$ws.Range("S8").Value = 0.01954992047285317
$ws.Range("T8").Value = 0.01954992047285318
$ws.Range("I9").Value = 0.3878984662564351
$ws.Range("J9").Value = 0.3878984662564351
$ws.Range("M9").Value = 81.23559966666666
$ws.Range("N9").Value = 243.706799
$ws.Range("O9").Value = 0.6800256617838641
$ws.Range("P9").Value = 0.6800256617838643
$ws.Range("Q9").Value = 73.13291724911433
$ws.Range("R9").Value = 658.196255242029
$ws.Range("S9").Value = 0.2637809112209782
$ws.Range("T9").Value = 0.2637809112209782
$ws.Range("G10").Value = 0.8083133333333334
$ws.Range("H10").Value = 2.42494
$ws.Range("I10").Value = 0.3482822152503414
$ws.Range("J10").Value = 0.3482822152503415
$ws.Range("M10").Value = 22.91402233333334
$ws.Range("N10").Value = 68.74206700000001
$ws.Range("O10").Value = 0.191813973987922
$ws.Range("P10").Value = 0.191813973987922
$ws.Range("Q10").Value = 18.52170977233111
$ws.Range("R10").Value = 166.69538795098
$ws.Range("S10").Value = 0.06680539577648485
$ws.Range("T10").Value = 0.06680539577648485
$ws.Range("G11").Value = 0.8083133333333334
$ws.Range("H11").Value = 2.42494
$ws.Range("I11").Value = 0.3482822152503414
$ws.Range("J11").Value = 0.3482822152503415
$ws.Range("O11").Value = 0.07776078244711707
$ws.Range("P11").Value = 0.07776078244711705
$ws.Range("Q11").Value = 7.508642953435557
$ws.Range("R11").Value = 67.57778658092001
$ws.Range("S11").Value = 0.0270826975702818
$ws.Range("T11").Value = 0.0270826975702818
$ws.Range("G12").Value = 0.8083133333333334
$ws.Range("H12").Value = 2.42494
$ws.Range("I12").Value = 0.3482822152503414
$ws.Range("J12").Value = 0.3482822152503415
$ws.Range("M12").Value = 6.020714333333333
$ws.Range("N12").Value = 18.062143
$ws.Range("O12").Value = 0.05039958178109668
$ws.Range("P12").Value = 0.05039958178109668
$ws.Range("Q12").Value = 4.866623671824445
$ws.Range("R12").Value = 43.79961304642001
$ws.Range("S12").Value = 0.0175532779904111
$ws.Range("T12").Value = 0.0175532779904111
$ws.Range("G13").Value = 0.8083133333333334
$ws.Range("H13").Value = 2.42494
$ws.Range("I13").Value = 0.3482822152503414
$ws.Range("J13").Value = 0.3482822152503415
$ws.Range("M13").Value = 81.23559966666666
$ws.Range("N13").Value = 243.706799
$ws.Range("O13").Value = 0.6800256617838641
$ws.Range("P13").Value = 0.6800256617838643
$ws.Range("Q13").Value = 65.66381835189556
$ws.Range("R13").Value = 590.97436516706
$ws.Range("S13").Value = 0.2368408439131637
$ws.Range("T13").Value = 0.2368408439131637
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.116974
$ws.Range("H14").Value = 0.350922
$ws.Range("I14").Value = 0.05040120231431718
$ws.Range("J14").Value = 0.05040120231431718
$ws.Range("M14").Value = 22.91402233333334
$ws.Range("N14").Value = 68.74206700000001
$ws.Range("O14").Value = 0.191813973987922
$ws.Range("P14").Value = 0.191813973987922
$ws.Range("Q14").Value = 2.680344848419334
$ws.Range("R14").Value = 24.123103635774
$ws.Range("S14").Value = 0.009667654909678431
$ws.Range("T14").Value = 0.009667654909678431
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.116974
$ws.Range("H15").Value = 0.350922
$ws.Range("I15").Value = 0.05040120231431718
$ws.Range("J15").Value = 0.05040120231431718
$ws.Range("O15").Value = 0.07776078244711707
$ws.Range("P15").Value = 0.07776078244711705
$ws.Range("Q15").Value = 1.086603380910667
$ws.Range("R15").Value = 9.779430428196001
$ws.Range("S15").Value = 0.003919236928236752
$ws.Range("T15").Value = 0.003919236928236751
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.116974
$ws.Range("H16").Value = 0.350922
$ws.Range("I16").Value = 0.05040120231431718
$ws.Range("J16").Value = 0.05040120231431718
$ws.Range("M16").Value = 6.020714333333333
$ws.Range("N16").Value = 18.062143
$ws.Range("O16").Value = 0.05039958178109668
$ws.Range("P16").Value = 0.05039958178109668
$ws.Range("Q16").Value = 0.7042670384273334
$ws.Range("R16").Value = 6.338403345845999
$ws.Range("S16").Value = 0.002540199517906028
$ws.Range("T16").Value = 0.002540199517906028
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.116974
$ws.Range("H17").Value = 0.350922
$ws.Range("I17").Value = 0.05040120231431718
$ws.Range("J17").Value = 0.05040120231431718
$ws.Range("M17").Value = 81.23559966666666
$ws.Range("N17").Value = 243.706799
$ws.Range("O17").Value = 0.6800256617838641
$ws.Range("P17").Value = 0.6800256617838643
$ws.Range("Q17").Value = 9.502453035408667
$ws.Range("R17").Value = 85.522077318678
$ws.Range("S17").Value = 0.03427411095849597
$ws.Range("T17").Value = 0.03427411095849597
